# Update profit data after running on 2025-09-15:
# append a new row (29) with the date and profit figure for that day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use text number format while assigning the value so Excel doesn't
# auto-convert the date-like string "09/15/2025" into a date serial
# number, then clear the temporary formatting so no extra style is
# left behind on the cell (matches how the other date cells in column
# A are stored as plain text).
$ws.Range("A29").NumberFormat = "@"
$ws.Range("A29").Value = "09/15/2025"
$ws.Range("A29").ClearFormats()

$ws.Range("B29").Value = 15297.89
